$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value cell (B8) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Include #0 sheet: update the "System URI" value cell (B4) ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R20-Pays/FHIR/TRE-R20-Pays"

# --- Include #1 sheet: update the "System URI" value cell (B4) ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R268-PaysProvenanceISO/FHIR/TRE-R268-PaysProvenanceISO"

# --- Include #2 sheet: update the "System URI" value cell (B4) ---
$wsInc2 = $wb.Worksheets.Item("Include #2")
$wsInc2.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R89-RegroupementPays/FHIR/TRE-R89-RegroupementPays"
